$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new wishlist entry as row 19 (Titulo / Autor / Editorial)
$ws.Range("A19").Value = "Nombres y animales"
$ws.Range("B19").Value = "Rita Indiana"
